{"js": "// Add quotations to column/table names (and a couple of small typo fixes\n// that were bundled into the same commit): a missing/extra hyphen-space\n// around \"Police(-| )deaths\" / \"gun(-| )death\", and a \"guns_deaths\" ->\n// \"gun_deaths\" typo fix.\n//\n// Strategy: the affected sentences are each a single, plain (no bold/\n// italic/hyperlink) paragraph, so we find each paragraph by a stable\n// substring of its original text and replace the *entire* paragraph text\n// in one shot with the fully-punctuated target text. This reproduces the\n// visible/textual result of the diff without depending on exactly how\n// Word happened to split the sentence into runs.\n\nconst LQ = \"\\u201c\"; // \u201c  left double quotation mark\nconst RQ = \"\\u201d\"; // \u201d  right double quotation mark\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfunction findParagraph(marker) {\n  const p = paragraphs.items.find((it) => it.text.indexOf(marker) !== -1);\n  if (!p) {\n    throw new Error(\"Could not find paragraph containing: \" + marker);\n  }\n  return p;\n}\n\nconst edits = [\n  {\n    marker: \"We used two different data sources from\",\n    text:\n      \"We used two different data sources from data.world, which are listed below. \" +\n      \"The Gun-deaths data source has death count by county and state dated from 1999 to 2016. \" +\n      \"Police-deaths data source has police deaths by name, date, cause of death, and by the police department dated from 1971 to 2016. \" +\n      \"The Gun-death data was in a .tsv file format, and Police department deaths in csv file format. \",\n  },\n  {\n    marker: \"The first step taken was separating the state\",\n    text:\n      \"The first step taken was separating the state from the county in the gun-death data source because \" +\n      \"the state is the unique identifier needed to join the two data sources. \",\n  },\n  {\n    marker: \"We removed columns in the police-deaths data sources\",\n    text:\n      \"We removed columns in the police-deaths data sources, columns removed are person, \" +\n      RQ + \"dept\" + RQ + \", \" + LQ + \"eow\" + RQ + \", \" + LQ + \"cause\" + RQ + \", \" + LQ + \"year\" + RQ + \", and \" + LQ + \"dept name\" + RQ + \". \" +\n      'In the gun-death data source, we replaced \"unreliable\" values with \"null\" on the \"crude_rate\" column. ' +\n      'Then replaced \"null\" with blanks because we encountered a problem when uploading to PostgreSQL. ',\n  },\n  {\n    marker: \"Using PostgreSQL, we created a new database\",\n    text:\n      'Using PostgreSQL, we created a new database, \"etl_database.\" We created 2 tables ' +\n      LQ + \"gun_deaths\" + RQ + \" table and \" + LQ + \"police_deaths\" + RQ + \" table with matching column headers to the csv files. \",\n  },\n  {\n    marker: \"columns are varchar data types\",\n    text:\n      LQ + \"county\" + RQ + \", \" + LQ + \"state\" + RQ + \", and \" + LQ + \"cause\" + RQ + \" columns are varchar data types. \" +\n      \"Population and deaths are integer data types. \" +\n      LQ + \"crude_rate\" + RQ + \" is a decimal data type, date column is a date data type, and \" + LQ + \"canine\" + RQ + \" a Boolean data type.\",\n  },\n  {\n    marker: \"We created two queries and assigned each of them an alias\",\n    text:\n      \"We created two queries and assigned each of them an alias. The first query was to get the count of each death type by state and assigned it to alias \" +\n      LQ + \"pd\" + RQ + \" on the \" + LQ + \"police_deaths\" + RQ + \" tables. \" +\n      \"The second query we selected all the count columns for type of death from the \" + LQ + \"police_deaths\" + RQ + \" table, \" +\n      \"then summed the deaths by state in the \" + LQ + \"gun_deaths\" + RQ + \" table and attached it to \" + LQ + \"nopd\" + RQ + \" alias, \" +\n      \"and joined the \" + LQ + \"pd\" + RQ + \" table to \" + LQ + \"nopd\" + RQ + \" table on column state. \" +\n      \"To finalize our database and to hold our query we create a view and named it \" + LQ + \"police_death_summary.\" + RQ,\n  },\n];\n\nfor (const edit of edits) {\n  const p = findParagraph(edit.marker);\n  p.getRange().insertText(edit.text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Add quotations to column/table names (and a couple of small typo fixes\n# that were bundled into the same commit): a missing/extra hyphen-space\n# around \"Police(-| )deaths\" / \"gun(-| )death\", and a \"guns_deaths\" ->\n# \"gun_deaths\" typo fix.\n#\n# We use Word's Find/Replace (wdReplaceAll = 2) scoped to the whole\n# document. Each search string is a long, unique substring of its\n# surrounding sentence so it can only ever match the one intended spot,\n# and the replacement substitutes in the fully-punctuated text (curly\n# quotes \"\\u201c\" / \"\\u201d\", matching Word's smart-quote AutoCorrect output).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text \"Police -deaths\" \"Police-deaths\"\n\nReplace-Text \"the gun- death data source\" \"the gun-death data source\"\n\nReplace-Text `\n    \"columns removed are person, dept, eow, cause, year, and dept name.\" `\n    \"columns removed are person, \u201ddept\u201d, \u201ceow\u201d, \u201ccause\u201d, \u201cyear\u201d, and \u201cdept name\u201d.\"\n\nReplace-Text `\n    \"We created 2 tables gun_deaths table and police_deaths table\" `\n    \"We created 2 tables \u201cgun_deaths\u201d table and \u201cpolice_deaths\u201d table\"\n\nReplace-Text `\n    \"County, state, and cause columns are varchar data types. Population and deaths are integer data types. Crude_rate is a decimal data type, date column is a date data type, and canine a Boolean data type.\" `\n    \"\u201ccounty\u201d, \u201cstate\u201d, and \u201ccause\u201d columns are varchar data types. Population and deaths are integer data types. \u201ccrude_rate\u201d is a decimal data type, date column is a date data type, and \u201ccanine\u201d a Boolean data type.\"\n\nReplace-Text `\n    \"assigned it to alias pd on the police_deaths tables.\" `\n    \"assigned it to alias \u201cpd\u201d on the \u201cpolice_deaths\u201d tables.\"\n\nReplace-Text `\n    \"type of death from the police_deaths table, then summed the deaths by state in the guns_deaths table and attached it to nopd alias, and joined the pd table to nopd table on column state.\" `\n    \"type of death from the \u201cpolice_deaths\u201d table, then summed the deaths by state in the \u201cgun_deaths\u201d table and attached it to \u201cnopd\u201d alias, and joined the \u201cpd\u201d table to \u201cnopd\u201d table on column state.\"\n\nReplace-Text `\n    \"we create a view and named it police_death_summary.\" `\n    \"we create a view and named it \u201cpolice_death_summary.\u201d\"\n"}
